{"js": "// Replace 100 arithmetic-problem cell texts per the diff mapping.\n// Each oldText is unique within the original document (verified offline),\n// so we resolve ALL 100 target ranges against the ORIGINAL, unmodified\n// document first, then apply the text replacements. Resolving every range\n// up front (rather than interleaving search+replace) avoids a later\n// search accidentally matching text that an earlier replacement just\n// inserted (e.g. old \"5+45=\" is a substring of new \"25+45=\").\nconst replacements = [[\"71+1=\", \"60-12=\"], [\"85-13=\", \"70-48=\"], [\"80-58=\", \"88-69=\"], [\"88-28=\", \"96-88=\"], [\"31-17=\", \"5+15=\"], [\"80+4=\", \"68-25=\"], [\"23+68=\", \"53+35=\"], [\"49+8=\", \"8+38=\"], [\"94-5=\", \"55-24=\"], [\"77-58=\", \"29+34=\"], [\"52+30=\", \"56-11=\"], [\"96-44=\", \"39-26=\"], [\"70-47=\", \"83-16=\"], [\"78-18=\", \"55-49=\"], [\"42-16=\", \"76-40=\"], [\"13+13=\", \"48-26=\"], [\"54-7=\", \"73-4=\"], [\"70-38=\", \"17+76=\"], [\"77-69=\", \"57+7=\"], [\"38+44=\", \"37+22=\"], [\"73-40=\", \"50+36=\"], [\"37+53=\", \"16-13=\"], [\"6+88=\", \"40+29=\"], [\"46+51=\", \"48+18=\"], [\"84+7=\", \"75-68=\"], [\"20-17=\", \"41-22=\"], [\"86-55=\", \"21+68=\"], [\"67-41=\", \"87+9=\"], [\"17+71=\", \"16+28=\"], [\"12+61=\", \"70+26=\"], [\"4+23=\", \"50+9=\"], [\"0+84=\", \"17-14=\"], [\"57+13=\", \"19+2=\"], [\"75-2=\", \"33+17=\"], [\"70-41=\", \"69-44=\"], [\"13-5=\", \"90-56=\"], [\"66-52=\", \"4+10=\"], [\"27-22=\", \"89-71=\"], [\"19+0=\", \"46-3=\"], [\"69+8=\", \"91+0=\"], [\"89-53=\", \"57+1=\"], [\"43+56=\", \"51+9=\"], [\"59-13=\", \"35-32=\"], [\"59-7=\", \"29+68=\"], [\"17+59=\", \"25+45=\"], [\"65-55=\", \"94-82=\"], [\"32+6=\", \"1+25=\"], [\"14+1=\", \"18+70=\"], [\"26+25=\", \"70+16=\"], [\"27+1=\", \"23+48=\"], [\"30+64=\", \"67-39=\"], [\"24+22=\", \"98-36=\"], [\"26+69=\", \"92-35=\"], [\"9+41=\", \"40-11=\"], [\"29+56=\", \"82-52=\"], [\"6+38=\", \"87-52=\"], [\"97-4=\", \"1+44=\"], [\"25-7=\", \"11+67=\"], [\"13+80=\", \"13+26=\"], [\"53+19=\", \"18+71=\"], [\"83-66=\", \"49+20=\"], [\"25+16=\", \"91-57=\"], [\"26+66=\", \"3+26=\"], [\"51+23=\", \"55+34=\"], [\"52+32=\", \"24-22=\"], [\"95-41=\", \"37+0=\"], [\"74-27=\", \"38+4=\"], [\"98-87=\", \"73-57=\"], [\"98-30=\", \"56-39=\"], [\"15+40=\", \"28+21=\"], [\"4+48=\", \"74-24=\"], [\"56+8=\", \"65-25=\"], [\"11+69=\", \"77-15=\"], [\"79+7=\", \"49-25=\"], [\"16+77=\", \"16+58=\"], [\"40-7=\", \"34+4=\"], [\"5+1=\", \"85-15=\"], [\"89-32=\", \"17+1=\"], [\"26+65=\", \"14+44=\"], [\"87-54=\", \"20+69=\"], [\"50+49=\", \"20+76=\"], [\"53+3=\", \"38+13=\"], [\"29+28=\", \"4+50=\"], [\"11+79=\", \"52+29=\"], [\"25+73=\", \"41+47=\"], [\"31+59=\", \"86-75=\"], [\"57-24=\", \"67+28=\"], [\"46+49=\", \"13-12=\"], [\"11+22=\", \"19+13=\"], [\"90-87=\", \"61+29=\"], [\"59-24=\", \"40+51=\"], [\"23+62=\", \"95-67=\"], [\"92-16=\", \"0+60=\"], [\"31-28=\", \"57+38=\"], [\"2+25=\", \"86-59=\"], [\"4+64=\", \"43+13=\"], [\"49-37=\", \"77+3=\"], [\"5+45=\", \"33+66=\"], [\"21+25=\", \"78-14=\"], [\"83-27=\", \"53+6=\"]];\n\nconst body = context.document.body;\nconst ranges = [];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(`Expected exactly 1 match for \"${oldText}\", got ${results.items.length}`);\n  }\n\n  ranges.push({ range: results.items[0], newText });\n}\n\nfor (const { range, newText } of ranges) {\n  range.insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-problem cell texts per the diff mapping.\n# The table is a fixed 20x5 grid (row-major) and every target cell holds\n# exactly one run of plain text, so we walk the table structurally by\n# (row, column) rather than doing a global text search-and-replace --\n# that sidesteps any accidental substring collision between an old value\n# still to be replaced and a new value already written (e.g. old \"5+45=\"\n# is a substring of new \"25+45=\").\n$oldValues = @(\n    \"71+1=\",\n    \"85-13=\",\n    \"80-58=\",\n    \"88-28=\",\n    \"31-17=\",\n    \"80+4=\",\n    \"23+68=\",\n    \"49+8=\",\n    \"94-5=\",\n    \"77-58=\",\n    \"52+30=\",\n    \"96-44=\",\n    \"70-47=\",\n    \"78-18=\",\n    \"42-16=\",\n    \"13+13=\",\n    \"54-7=\",\n    \"70-38=\",\n    \"77-69=\",\n    \"38+44=\",\n    \"73-40=\",\n    \"37+53=\",\n    \"6+88=\",\n    \"46+51=\",\n    \"84+7=\",\n    \"20-17=\",\n    \"86-55=\",\n    \"67-41=\",\n    \"17+71=\",\n    \"12+61=\",\n    \"4+23=\",\n    \"0+84=\",\n    \"57+13=\",\n    \"75-2=\",\n    \"70-41=\",\n    \"13-5=\",\n    \"66-52=\",\n    \"27-22=\",\n    \"19+0=\",\n    \"69+8=\",\n    \"89-53=\",\n    \"43+56=\",\n    \"59-13=\",\n    \"59-7=\",\n    \"17+59=\",\n    \"65-55=\",\n    \"32+6=\",\n    \"14+1=\",\n    \"26+25=\",\n    \"27+1=\",\n    \"30+64=\",\n    \"24+22=\",\n    \"26+69=\",\n    \"9+41=\",\n    \"29+56=\",\n    \"6+38=\",\n    \"97-4=\",\n    \"25-7=\",\n    \"13+80=\",\n    \"53+19=\",\n    \"83-66=\",\n    \"25+16=\",\n    \"26+66=\",\n    \"51+23=\",\n    \"52+32=\",\n    \"95-41=\",\n    \"74-27=\",\n    \"98-87=\",\n    \"98-30=\",\n    \"15+40=\",\n    \"4+48=\",\n    \"56+8=\",\n    \"11+69=\",\n    \"79+7=\",\n    \"16+77=\",\n    \"40-7=\",\n    \"5+1=\",\n    \"89-32=\",\n    \"26+65=\",\n    \"87-54=\",\n    \"50+49=\",\n    \"53+3=\",\n    \"29+28=\",\n    \"11+79=\",\n    \"25+73=\",\n    \"31+59=\",\n    \"57-24=\",\n    \"46+49=\",\n    \"11+22=\",\n    \"90-87=\",\n    \"59-24=\",\n    \"23+62=\",\n    \"92-16=\",\n    \"31-28=\",\n    \"2+25=\",\n    \"4+64=\",\n    \"49-37=\",\n    \"5+45=\",\n    \"21+25=\",\n    \"83-27=\"\n)\n\n$newValues = @(\n    \"60-12=\",\n    \"70-48=\",\n    \"88-69=\",\n    \"96-88=\",\n    \"5+15=\",\n    \"68-25=\",\n    \"53+35=\",\n    \"8+38=\",\n    \"55-24=\",\n    \"29+34=\",\n    \"56-11=\",\n    \"39-26=\",\n    \"83-16=\",\n    \"55-49=\",\n    \"76-40=\",\n    \"48-26=\",\n    \"73-4=\",\n    \"17+76=\",\n    \"57+7=\",\n    \"37+22=\",\n    \"50+36=\",\n    \"16-13=\",\n    \"40+29=\",\n    \"48+18=\",\n    \"75-68=\",\n    \"41-22=\",\n    \"21+68=\",\n    \"87+9=\",\n    \"16+28=\",\n    \"70+26=\",\n    \"50+9=\",\n    \"17-14=\",\n    \"19+2=\",\n    \"33+17=\",\n    \"69-44=\",\n    \"90-56=\",\n    \"4+10=\",\n    \"89-71=\",\n    \"46-3=\",\n    \"91+0=\",\n    \"57+1=\",\n    \"51+9=\",\n    \"35-32=\",\n    \"29+68=\",\n    \"25+45=\",\n    \"94-82=\",\n    \"1+25=\",\n    \"18+70=\",\n    \"70+16=\",\n    \"23+48=\",\n    \"67-39=\",\n    \"98-36=\",\n    \"92-35=\",\n    \"40-11=\",\n    \"82-52=\",\n    \"87-52=\",\n    \"1+44=\",\n    \"11+67=\",\n    \"13+26=\",\n    \"18+71=\",\n    \"49+20=\",\n    \"91-57=\",\n    \"3+26=\",\n    \"55+34=\",\n    \"24-22=\",\n    \"37+0=\",\n    \"38+4=\",\n    \"73-57=\",\n    \"56-39=\",\n    \"28+21=\",\n    \"74-24=\",\n    \"65-25=\",\n    \"77-15=\",\n    \"49-25=\",\n    \"16+58=\",\n    \"34+4=\",\n    \"85-15=\",\n    \"17+1=\",\n    \"14+44=\",\n    \"20+69=\",\n    \"20+76=\",\n    \"38+13=\",\n    \"4+50=\",\n    \"52+29=\",\n    \"41+47=\",\n    \"86-75=\",\n    \"67+28=\",\n    \"13-12=\",\n    \"19+13=\",\n    \"61+29=\",\n    \"40+51=\",\n    \"95-67=\",\n    \"0+60=\",\n    \"57+38=\",\n    \"86-59=\",\n    \"43+13=\",\n    \"77+3=\",\n    \"33+66=\",\n    \"78-14=\",\n    \"53+6=\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $t.Cell($r, $c)\n        # Cell.Range.Text carries Word's trailing end-of-cell mark (CR + BEL,\n        # chars 13/7) after the visible text -- strip it before comparing.\n        $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n        $expectedOld = $oldValues[$i]\n        if ($current -ne $expectedOld) {\n            throw \"Cell ($r,$c) text mismatch: expected '$expectedOld' but found '$current'\"\n        }\n        $cell.Range.Text = $newValues[$i]\n        $i = $i + 1\n    }\n}\n\nWrite-Output \"Replaced $i cells\"\n"}
